# semana 22 de 2025
# Update the Esperado (C), Observado (D) and valor p (E) columns on Sheet1
# to reflect the latest epidemiological figures for week 22 of 2025.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 3).Value = 2
$ws.Cells.Item(2, 5).Value = 0.09

$ws.Cells.Item(3, 3).Value = 1
$ws.Cells.Item(3, 5).Value = 0.37

$ws.Cells.Item(4, 3).Value = 9
$ws.Cells.Item(4, 5).Value = 0.12

$ws.Cells.Item(5, 3).Value = 2
$ws.Cells.Item(5, 4).Value = 14

$ws.Cells.Item(6, 3).Value = 2
$ws.Cells.Item(6, 5).Value = 0.18

$ws.Cells.Item(7, 4).Value = 1
$ws.Cells.Item(7, 5).Value = 0

$ws.Cells.Item(8, 3).Value = 41
$ws.Cells.Item(8, 4).Value = 38
$ws.Cells.Item(8, 5).Value = 0.06

$ws.Cells.Item(9, 3).Value = 0
$ws.Cells.Item(9, 5).Value = 1

$ws.Cells.Item(10, 4).Value = 3
$ws.Cells.Item(10, 5).Value = 0.06

$ws.Cells.Item(11, 4).Value = 4
$ws.Cells.Item(11, 5).Value = 0.17

$ws.Cells.Item(12, 3).Value = 12
$ws.Cells.Item(12, 4).Value = 4
$ws.Cells.Item(12, 5).Value = 0.01

$ws.Cells.Item(14, 4).Value = 1

$ws.Cells.Item(15, 3).Value = 0
$ws.Cells.Item(15, 5).Value = 0

$ws.Cells.Item(16, 3).Value = 9
$ws.Cells.Item(16, 4).Value = 8
$ws.Cells.Item(16, 5).Value = 0.13

$ws.Cells.Item(18, 3).Value = 5
$ws.Cells.Item(18, 4).Value = 3
$ws.Cells.Item(18, 5).Value = 0.14

$ws.Cells.Item(19, 4).Value = 0
$ws.Cells.Item(19, 5).Value = 0.37

$ws.Cells.Item(20, 4).Value = 2

$ws.Cells.Item(22, 3).Value = 5
$ws.Cells.Item(22, 4).Value = 1
$ws.Cells.Item(22, 5).Value = 0.03

$ws.Cells.Item(23, 3).Value = 1
$ws.Cells.Item(23, 5).Value = 0.37

$ws.Cells.Item(28, 3).Value = 0
$ws.Cells.Item(28, 5).Value = 1

$ws.Cells.Item(29, 3).Value = 1
$ws.Cells.Item(29, 4).Value = 1
$ws.Cells.Item(29, 5).Value = 0.37

$ws.Cells.Item(30, 4).Value = 1
$ws.Cells.Item(30, 5).Value = 0

$ws.Cells.Item(31, 3).Value = 7
$ws.Cells.Item(31, 5).Value = 0.02

$ws.Cells.Item(32, 3).Value = 10
$ws.Cells.Item(32, 4).Value = 4
$ws.Cells.Item(32, 5).Value = 0.02

$ws.Cells.Item(33, 3).Value = 8
$ws.Cells.Item(33, 4).Value = 9
$ws.Cells.Item(33, 5).Value = 0.12
